# "Criada a classe gráfico"
#
# Statistics row (row 18) in the "TabelaHMNaoNP" table gets its extra
# descriptive-statistics formulas reshuffled: a new "Moda" (MODE.SNGL)
# computation is introduced in column E, the old column-D formula
# (AVEDEV) is replaced by MEDIAN, and the STDEV.P / VAR.P computations
# that used to live in E/F shift one column to the right (into F/G),
# displacing the old "Coeficiente de Variação" (E18/C18) formula that
# used to sit in G18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D: Desvio Médio -> now computed with MEDIAN instead of AVEDEV
$ws.Range("D18").Formula = "=MEDIAN(B1:B18)"

# Column E: new "Moda" calculated-column formula (MODE.SNGL)
$ws.Range("E18").Formula = "=MODE.SNGL(B1:B18)"

# Column F: now holds the old "Desvio Padrão" (STDEV.P) formula
$ws.Range("F18").Formula = "=STDEV.P(B1:B18)"

# Column G: now holds the old "Variância Populacional" (VAR.P) formula,
# replacing the previous "Coeficiente de Variação" (E18/C18) formula
$ws.Range("G18").Formula = "=VAR.P(B1:B18)"

# H18 (Variância Amostral) is untouched.

# The table treats E18 as a one-off "calculated column" formula (it is
# the sole populated cell in that table column) — mirror that, and mark
# the resulting "inconsistent formula" warning as ignored, same as a
# user dismissing the green error triangle in the real UI.
$lo = $ws.ListObjects.Item(1)
try {
    $col5 = $lo.ListColumns.Item(5)
    $col5.Formula = "=MODE.SNGL(B1:B18)"
} catch {
}

try {
    $err = $ws.Range("E18").Errors.Item(4)
    $err.Ignore = $true
} catch {
}

# Move the active selection to D19, matching where the user clicked
# next after editing the row.
$ws.Range("D19").Select() | Out-Null
